# Rerun with new APC variants: refresh the per-histology p-values
# (column C, "pvalue") computed for age/sex/race/ancestry/ethnicity/
# OS_years/EFS_years/tmb across the affected histology tabs.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Cells.Item(2, 3).Value = 0.0896185966222548
$ws.Cells.Item(3, 3).Value = 0.0208619595159617
$ws.Cells.Item(4, 3).Value = 0.90258835020493
$ws.Cells.Item(5, 3).Value = 0.617880140186009
$ws.Cells.Item(6, 3).Value = 0.276767428792189
$ws.Cells.Item(7, 3).Value = 0.732449677492432
$ws.Cells.Item(8, 3).Value = 0.282133084993859
$ws.Cells.Item(9, 3).Value = 0.284369753577153

$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Cells.Item(2, 3).Value = 0.533238741934394
$ws.Cells.Item(3, 3).Value = 0.408212560386474
$ws.Cells.Item(4, 3).Value = 0.282003986351813
$ws.Cells.Item(5, 3).Value = 0.615904868078785
$ws.Cells.Item(6, 3).Value = 0.231809736157562
$ws.Cells.Item(9, 3).Value = 0.353287117317915

$ws = $wb.Worksheets.Item("Other tumor")
$ws.Cells.Item(2, 3).Value = 0.0907981961597089
$ws.Cells.Item(4, 3).Value = 0.426613045834412
$ws.Cells.Item(5, 3).Value = 0.81427569093866
$ws.Cells.Item(6, 3).Value = 0.780633752602614
$ws.Cells.Item(7, 3).Value = 0.878298596491423
$ws.Cells.Item(8, 3).Value = 0.759361447436373
$ws.Cells.Item(9, 3).Value = 0.739996268959378

$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Cells.Item(2, 3).Value = 0.713705357283287
$ws.Cells.Item(4, 3).Value = 0.987141116454481
$ws.Cells.Item(5, 3).Value = 0.946801804848752
$ws.Cells.Item(6, 3).Value = 0.639519327674786
$ws.Cells.Item(7, 3).Value = 0.00369227398780271
$ws.Cells.Item(8, 3).Value = 0.0034650073118416
$ws.Cells.Item(9, 3).Value = 0.0878796824057952

$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Cells.Item(2, 3).Value = 0.498965426816437
$ws.Cells.Item(3, 3).Value = 0.0565098418868229
$ws.Cells.Item(4, 3).Value = 0.999999999999979
$ws.Cells.Item(5, 3).Value = 0.631203160659328
$ws.Cells.Item(6, 3).Value = 0.999999999999973
$ws.Cells.Item(7, 3).Value = 0.23832192708036
$ws.Cells.Item(8, 3).Value = 0.419110881762562
$ws.Cells.Item(9, 3).Value = 0.167425229268679

$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Cells.Item(2, 3).Value = 0.233152960277338
$ws.Cells.Item(4, 3).Value = 0.576531536291638
$ws.Cells.Item(5, 3).Value = 0.368336454051302
$ws.Cells.Item(6, 3).Value = 0.8325789171934
$ws.Cells.Item(7, 3).Value = 0.154349260156668
$ws.Cells.Item(8, 3).Value = 0.144877451159352
$ws.Cells.Item(9, 3).Value = 0.147476037454092

$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Cells.Item(2, 3).Value = 0.0625772201692652
$ws.Cells.Item(3, 3).Value = 0.116505436478945
$ws.Cells.Item(4, 3).Value = 0.35325139146675
$ws.Cells.Item(5, 3).Value = 0.575562395767317
$ws.Cells.Item(6, 3).Value = 0.402895113668335
$ws.Cells.Item(7, 3).Value = 0.516918809678844
$ws.Cells.Item(8, 3).Value = 0.886473974807497
$ws.Cells.Item(9, 3).Value = 0.201167318853901

$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Cells.Item(6, 3).Value = 0.999999999999972

Write-Host "Updated p-values on 8 histology sheets (52 cells)."
